# Update "想去人数" (want-to-go count) figures for two events that were
# refreshed when the site data was regenerated (gh-pages output at 456a3b4).
#
# 展览 (Exhibitions) sheet: row 3 is the
#   "合肥·第九届环形宇宙动漫游戏嘉年华" exhibition -> F3 2525->2538
# 演出 (Shows) sheet: row 2 is the
#   "四月是你的谎言" concert -> F2 114->115
# 全部类型 (All types) sheet aggregates both of the above, so its matching
# rows need the same updates: F3 114->115, F7 2525->2538

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 2538

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 115

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 115
$wsAll.Range("F7").Value = 2538
